$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 22, pushing existing rows 22-39 down to 23-40
$ws.Rows.Item(22).Insert()

# Populate the newly inserted row 22 with the new record.
# Columns A,B,C,E,F,G,I,R carry the same constant values as the rest of this
# dataset (Agrícola del Norte S.A. de Arica / Poroto verde / Primera / Hortaliza).
$ws.Cells.Item(22, 1).Value = 1
$ws.Cells.Item(22, 2).Value = "Agrícola del Norte S.A. de Arica"
$ws.Cells.Item(22, 3).Value = "Arica y Parinacota"
$ws.Cells.Item(22, 4).Value = 44483
$ws.Cells.Item(22, 5).Value = 15
$ws.Cells.Item(22, 6).Value = 100112031
$ws.Cells.Item(22, 7).Value = "Poroto verde"
$ws.Cells.Item(22, 8).Value = "Sin especificar"
$ws.Cells.Item(22, 9).Value = "Primera"
$ws.Cells.Item(22, 10).Value = 1100
$ws.Cells.Item(22, 11).Value = 1200
$ws.Cells.Item(22, 12).Value = 1300
$ws.Cells.Item(22, 13).Value = 1250
$ws.Cells.Item(22, 14).Value = "$/kilo"
$ws.Cells.Item(22, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(22, 16).Value = 1250
$ws.Cells.Item(22, 17).Value = 1
$ws.Cells.Item(22, 18).Value = "Hortaliza"
